$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 95; $r++) {
    if ($r -eq 36) { continue }
    $cell = $ws.Cells.Item($r, 5)
    $current = $cell.Value2
    $cell.Value = $current - 1
}
